$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range('K2').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743717012329/Ggg0CQFNgb-camera-video-segment-driver-1743717014829.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLHAFRQHL%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T190016Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFVsYb1pK%2FPZjGbw6%2BCTlRsyY9p23pjlnu9Ycafk3TsuAiEAzKG7wFiGe%2Fy7UTC9Pylw2xp8oJEX%2FlfT4GKCEXz6Nbsq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDBqrq2agNGfdGcefNiq6A1WvOi07cMST%2FIyRPJkK3hVP%2BQmk2cEl1lKpU4DXDo4DifQI%2FJtyHFjTK4RzIvvoHS9IcMOSCSeiXsLkdkZIs9hZHFVgB%2Bq1ACdZ6Ijt8Y5h%2BhcBU3bnLRrBWzlAN2ucbnyPcJ4Wb8xNCvgSvbqKuf7qwkriR5aWEmDQt%2B5UsZ56V5uDi1j%2F%2F1x%2Bljl5Llpi%2FgIuNRwkgeP5Gl43XuR6IP9%2F664QsL9ADXt%2FGip2yjUY8zoV4bt1Oytfsupc3aRk8GP58jXJ7SsvN1gAti8w5QL1fPGZ19Pi84KUPu72Cf%2FRub54ctQl0RbrkoTLrmi7gQOB%2FPLMUdjkp3wVjwF5l4nA4CMMccLd01XT25QEjh%2B5PAxuiLya0DfXeDxDVL6JUT4%2F1FdEcmxV6tUACVsX5uFgYLEuae4YvWrDKw1IMUi%2BbK5ffv42pv5ZJ%2F8G034Oqedv0DIQXOVhEMJtFUi4SVVl6kUnl6E4A%2BDks8Dko%2BaEouzmTbSxhfjaSTG82IOlBuL1fgnS5NyP9aCBPW6Spr2WHP94zzoW6SflNq4dDNTgNQDvXCuUU7oKUXjBbWeJ8u1VPV5mWeHPMVow98LAvwY6pQEtz3oNp9eZI7Due%2FQQWFan4Fprxim3kzZdDJXuQNmFkMUoFXnNH2EnNIqjyj3oV%2BWZUVOKXexNP8x9LYmvAdbE6Unju%2BE3a0y773dji%2FWBmCJNsImUO%2F0XGRyyIcCoMijx2ST5rKJoRuKBRnW72VXK7UEYO%2BSiRApw%2B3notP4163zqsrPusLBQmfCDVe9MRKB0XNbOYoYoHLtsH0945Bsjg7TklRc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2003%3A00%3A16%20GMT&X-Amz-Signature=8a891839d78bafd202c22adb7ab61a46c64452199b74ff8d5d110f4ead52e2f7'
$ws.Range('K3').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743713574485/FgwSeOefdQ-camera-video-segment-driver-1743713576985.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLHAFRQHL%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T190016Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFVsYb1pK%2FPZjGbw6%2BCTlRsyY9p23pjlnu9Ycafk3TsuAiEAzKG7wFiGe%2Fy7UTC9Pylw2xp8oJEX%2FlfT4GKCEXz6Nbsq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDBqrq2agNGfdGcefNiq6A1WvOi07cMST%2FIyRPJkK3hVP%2BQmk2cEl1lKpU4DXDo4DifQI%2FJtyHFjTK4RzIvvoHS9IcMOSCSeiXsLkdkZIs9hZHFVgB%2Bq1ACdZ6Ijt8Y5h%2BhcBU3bnLRrBWzlAN2ucbnyPcJ4Wb8xNCvgSvbqKuf7qwkriR5aWEmDQt%2B5UsZ56V5uDi1j%2F%2F1x%2Bljl5Llpi%2FgIuNRwkgeP5Gl43XuR6IP9%2F664QsL9ADXt%2FGip2yjUY8zoV4bt1Oytfsupc3aRk8GP58jXJ7SsvN1gAti8w5QL1fPGZ19Pi84KUPu72Cf%2FRub54ctQl0RbrkoTLrmi7gQOB%2FPLMUdjkp3wVjwF5l4nA4CMMccLd01XT25QEjh%2B5PAxuiLya0DfXeDxDVL6JUT4%2F1FdEcmxV6tUACVsX5uFgYLEuae4YvWrDKw1IMUi%2BbK5ffv42pv5ZJ%2F8G034Oqedv0DIQXOVhEMJtFUi4SVVl6kUnl6E4A%2BDks8Dko%2BaEouzmTbSxhfjaSTG82IOlBuL1fgnS5NyP9aCBPW6Spr2WHP94zzoW6SflNq4dDNTgNQDvXCuUU7oKUXjBbWeJ8u1VPV5mWeHPMVow98LAvwY6pQEtz3oNp9eZI7Due%2FQQWFan4Fprxim3kzZdDJXuQNmFkMUoFXnNH2EnNIqjyj3oV%2BWZUVOKXexNP8x9LYmvAdbE6Unju%2BE3a0y773dji%2FWBmCJNsImUO%2F0XGRyyIcCoMijx2ST5rKJoRuKBRnW72VXK7UEYO%2BSiRApw%2B3notP4163zqsrPusLBQmfCDVe9MRKB0XNbOYoYoHLtsH0945Bsjg7TklRc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2003%3A00%3A16%20GMT&X-Amz-Signature=f5969b323c347a8e1429b496350f6003f092f76a5481e5a36e82dc8ff7ce2601'
$ws.Range('K4').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712999664/UzzO0YQTJX-camera-video-segment-driver-1743713002164.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLHAFRQHL%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T190016Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFVsYb1pK%2FPZjGbw6%2BCTlRsyY9p23pjlnu9Ycafk3TsuAiEAzKG7wFiGe%2Fy7UTC9Pylw2xp8oJEX%2FlfT4GKCEXz6Nbsq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDBqrq2agNGfdGcefNiq6A1WvOi07cMST%2FIyRPJkK3hVP%2BQmk2cEl1lKpU4DXDo4DifQI%2FJtyHFjTK4RzIvvoHS9IcMOSCSeiXsLkdkZIs9hZHFVgB%2Bq1ACdZ6Ijt8Y5h%2BhcBU3bnLRrBWzlAN2ucbnyPcJ4Wb8xNCvgSvbqKuf7qwkriR5aWEmDQt%2B5UsZ56V5uDi1j%2F%2F1x%2Bljl5Llpi%2FgIuNRwkgeP5Gl43XuR6IP9%2F664QsL9ADXt%2FGip2yjUY8zoV4bt1Oytfsupc3aRk8GP58jXJ7SsvN1gAti8w5QL1fPGZ19Pi84KUPu72Cf%2FRub54ctQl0RbrkoTLrmi7gQOB%2FPLMUdjkp3wVjwF5l4nA4CMMccLd01XT25QEjh%2B5PAxuiLya0DfXeDxDVL6JUT4%2F1FdEcmxV6tUACVsX5uFgYLEuae4YvWrDKw1IMUi%2BbK5ffv42pv5ZJ%2F8G034Oqedv0DIQXOVhEMJtFUi4SVVl6kUnl6E4A%2BDks8Dko%2BaEouzmTbSxhfjaSTG82IOlBuL1fgnS5NyP9aCBPW6Spr2WHP94zzoW6SflNq4dDNTgNQDvXCuUU7oKUXjBbWeJ8u1VPV5mWeHPMVow98LAvwY6pQEtz3oNp9eZI7Due%2FQQWFan4Fprxim3kzZdDJXuQNmFkMUoFXnNH2EnNIqjyj3oV%2BWZUVOKXexNP8x9LYmvAdbE6Unju%2BE3a0y773dji%2FWBmCJNsImUO%2F0XGRyyIcCoMijx2ST5rKJoRuKBRnW72VXK7UEYO%2BSiRApw%2B3notP4163zqsrPusLBQmfCDVe9MRKB0XNbOYoYoHLtsH0945Bsjg7TklRc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2003%3A00%3A16%20GMT&X-Amz-Signature=6ea43ab5fa524b015177a5beac29444b3ef4e8e67f4feb0e75d62c3f376bca5c'
$ws.Range('K5').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712546260/Jh2p7rn04n-camera-video-segment-driver-1743712548760.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLHAFRQHL%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T190016Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFVsYb1pK%2FPZjGbw6%2BCTlRsyY9p23pjlnu9Ycafk3TsuAiEAzKG7wFiGe%2Fy7UTC9Pylw2xp8oJEX%2FlfT4GKCEXz6Nbsq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDBqrq2agNGfdGcefNiq6A1WvOi07cMST%2FIyRPJkK3hVP%2BQmk2cEl1lKpU4DXDo4DifQI%2FJtyHFjTK4RzIvvoHS9IcMOSCSeiXsLkdkZIs9hZHFVgB%2Bq1ACdZ6Ijt8Y5h%2BhcBU3bnLRrBWzlAN2ucbnyPcJ4Wb8xNCvgSvbqKuf7qwkriR5aWEmDQt%2B5UsZ56V5uDi1j%2F%2F1x%2Bljl5Llpi%2FgIuNRwkgeP5Gl43XuR6IP9%2F664QsL9ADXt%2FGip2yjUY8zoV4bt1Oytfsupc3aRk8GP58jXJ7SsvN1gAti8w5QL1fPGZ19Pi84KUPu72Cf%2FRub54ctQl0RbrkoTLrmi7gQOB%2FPLMUdjkp3wVjwF5l4nA4CMMccLd01XT25QEjh%2B5PAxuiLya0DfXeDxDVL6JUT4%2F1FdEcmxV6tUACVsX5uFgYLEuae4YvWrDKw1IMUi%2BbK5ffv42pv5ZJ%2F8G034Oqedv0DIQXOVhEMJtFUi4SVVl6kUnl6E4A%2BDks8Dko%2BaEouzmTbSxhfjaSTG82IOlBuL1fgnS5NyP9aCBPW6Spr2WHP94zzoW6SflNq4dDNTgNQDvXCuUU7oKUXjBbWeJ8u1VPV5mWeHPMVow98LAvwY6pQEtz3oNp9eZI7Due%2FQQWFan4Fprxim3kzZdDJXuQNmFkMUoFXnNH2EnNIqjyj3oV%2BWZUVOKXexNP8x9LYmvAdbE6Unju%2BE3a0y773dji%2FWBmCJNsImUO%2F0XGRyyIcCoMijx2ST5rKJoRuKBRnW72VXK7UEYO%2BSiRApw%2B3notP4163zqsrPusLBQmfCDVe9MRKB0XNbOYoYoHLtsH0945Bsjg7TklRc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2003%3A00%3A16%20GMT&X-Amz-Signature=10a83c6dc85a963b7d79bd48d9585683e530188f06be4c276affb2175ab31f8f'
$ws.Range('K6').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743703000742/spTbAbfADq-camera-video-segment-driver-1743703003242.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLHAFRQHL%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T190016Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFVsYb1pK%2FPZjGbw6%2BCTlRsyY9p23pjlnu9Ycafk3TsuAiEAzKG7wFiGe%2Fy7UTC9Pylw2xp8oJEX%2FlfT4GKCEXz6Nbsq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDBqrq2agNGfdGcefNiq6A1WvOi07cMST%2FIyRPJkK3hVP%2BQmk2cEl1lKpU4DXDo4DifQI%2FJtyHFjTK4RzIvvoHS9IcMOSCSeiXsLkdkZIs9hZHFVgB%2Bq1ACdZ6Ijt8Y5h%2BhcBU3bnLRrBWzlAN2ucbnyPcJ4Wb8xNCvgSvbqKuf7qwkriR5aWEmDQt%2B5UsZ56V5uDi1j%2F%2F1x%2Bljl5Llpi%2FgIuNRwkgeP5Gl43XuR6IP9%2F664QsL9ADXt%2FGip2yjUY8zoV4bt1Oytfsupc3aRk8GP58jXJ7SsvN1gAti8w5QL1fPGZ19Pi84KUPu72Cf%2FRub54ctQl0RbrkoTLrmi7gQOB%2FPLMUdjkp3wVjwF5l4nA4CMMccLd01XT25QEjh%2B5PAxuiLya0DfXeDxDVL6JUT4%2F1FdEcmxV6tUACVsX5uFgYLEuae4YvWrDKw1IMUi%2BbK5ffv42pv5ZJ%2F8G034Oqedv0DIQXOVhEMJtFUi4SVVl6kUnl6E4A%2BDks8Dko%2BaEouzmTbSxhfjaSTG82IOlBuL1fgnS5NyP9aCBPW6Spr2WHP94zzoW6SflNq4dDNTgNQDvXCuUU7oKUXjBbWeJ8u1VPV5mWeHPMVow98LAvwY6pQEtz3oNp9eZI7Due%2FQQWFan4Fprxim3kzZdDJXuQNmFkMUoFXnNH2EnNIqjyj3oV%2BWZUVOKXexNP8x9LYmvAdbE6Unju%2BE3a0y773dji%2FWBmCJNsImUO%2F0XGRyyIcCoMijx2ST5rKJoRuKBRnW72VXK7UEYO%2BSiRApw%2B3notP4163zqsrPusLBQmfCDVe9MRKB0XNbOYoYoHLtsH0945Bsjg7TklRc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2003%3A00%3A16%20GMT&X-Amz-Signature=c8799707b58f9e7b73461231aa38685d5f708ac30c73518c29312070d6695861'
